# Fix formatting for the right lever data single-cell string (and the
# related C40:C46 "R Minutes" column that had been filled with the wrong
# values - actually the L-side minutes - by mistake).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected "R Minutes" values for rows 40-46 (column C)
$ws.Range("C40").Value = 5.507
$ws.Range("C41").Value = 8.074
$ws.Range("C42").Value = 11.705
$ws.Range("C43").Value = 13.469
$ws.Range("C44").Value = 99.842
$ws.Range("C45").Value = 99.859
$ws.Range("C46").Value = 113.329

# The "L ~ Minutes" single-cell string (A56) had erroneously contained a
# stale/incorrect list. Replace it with the correct L Minutes values
# (matching column C rows 4-38).
$ws.Range("A56").Value = "0.283, 0.883, 2.026, 2.696, 3.609, 4.574, 5.281, 5.696, 5.737, 5.905, 6.474, 7.157, 7.268, 7.346, 7.444, 8.133, 10.831, 11.467, 12.632, 12.735, 13.417, 15.738, 53.07, 67.125, 68.016, 69.649, 96.018, 96.553, 97.113, 98.248, 98.798, 99.974, 110.653, 111.559, 112.346"

# The "R ~ Minutes" single-cell string (A58) is updated to reflect the
# corrected R Minutes values above.
$ws.Range("A58").Value = "5.507, 8.074, 11.705, 13.469, 99.842, 99.859, 113.329"
